# Apply the "Updated cryptos list" GitHub Actions refresh to Sheet1.
# Price (col D) and Volume(1h) (col E) values are refreshed per row;
# rows 48/49 additionally swap which coin (Filecoin/ARBITRUM) occupies them.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '''69.596.92'
$ws.Range("E2").Value = '  -3.75%  '

# Row 3
$ws.Range("D3").Value = '''2.509.26'
$ws.Range("E3").Value = '  -5.18%  '

# Row 4
$ws.Range("E4").Value = '  -0.05%  '

# Row 5
$ws.Range("D5").Value = '''576.98'
$ws.Range("E5").Value = '  -2.41%  '

# Row 6
$ws.Range("D6").Value = '''167.20'
$ws.Range("E6").Value = '  -4.78%  '

# Row 7
$ws.Range("E7").Value = '  +0.05%  '

# Row 8
$ws.Range("D8").Value = '''0.515'
$ws.Range("E8").Value = '  -1.12%  '

# Row 9
$ws.Range("D9").Value = '''2.508.37'
$ws.Range("E9").Value = '  -5.18%  '

# Row 10
$ws.Range("E10").Value = '  -7.08%  '

# Row 11
$ws.Range("E11").Value = '  -0.71%  '

# Row 12
$ws.Range("E12").Value = '  -3.92%  '

# Row 13
$ws.Range("E13").Value = '  -1.98%  '

# Row 14
$ws.Range("D14").Value = '''2.982.85'
$ws.Range("E14").Value = '  -4.87%  '

# Row 15
$ws.Range("D15").Value = '''69.423.62'
$ws.Range("E15").Value = '  -3.89%  '

# Row 16
$ws.Range("E16").Value = '  -6.00%  '

# Row 17
$ws.Range("D17").Value = '''24.97'
$ws.Range("E17").Value = '  -3.99%  '

# Row 18
$ws.Range("D18").Value = '''2.518.89'
$ws.Range("E18").Value = '  -3.44%  '

# Row 19
$ws.Range("E19").Value = '  -6.58%  '

# Row 20
$ws.Range("E20").Value = '  -2.69%  '

# Row 21
$ws.Range("D21").Value = '''350.49'
$ws.Range("E21").Value = '  -5.32%  '

# Row 22
$ws.Range("E22").Value = '  -4.94%  '

# Row 23
$ws.Range("D23").Value = '''1.97'
$ws.Range("E23").Value = '  -4.76%  '

# Row 24
$ws.Range("E24").Value = '  -0.04%  '

# Row 25
$ws.Range("D25").Value = '''68.89'
$ws.Range("E25").Value = '  -3.59%  '

# Row 26
$ws.Range("E26").Value = '  -6.12%  '

# Row 27
$ws.Range("D27").Value = '''9.01'
$ws.Range("E27").Value = '  -7.12%  '

# Row 28
$ws.Range("D28").Value = '''2.638.66'
$ws.Range("E28").Value = '  -5.26%  '

# Row 29
$ws.Range("D29").Value = '''0.998'
$ws.Range("E29").Value = '  -0.32%  '

# Row 30
$ws.Range("D30").Value = '''0.0₃0903'
$ws.Range("E30").Value = '  -6.00%  '

# Row 31
$ws.Range("D31").Value = '''7.92'
$ws.Range("E31").Value = '  -1.54%  '

# Row 32
$ws.Range("D32").Value = '''478.58'
$ws.Range("E32").Value = '  -4.56%  '

# Row 33
$ws.Range("D33").Value = '''1.29'
$ws.Range("E33").Value = '  -0.26%  '

# Row 34
$ws.Range("E34").Value = '  -3.28%  '

# Row 35
$ws.Range("E35").Value = '  -0.04%  '

# Row 36
$ws.Range("E36").Value = '  -0.74%  '

# Row 37
$ws.Range("D37").Value = '''154.97'
$ws.Range("E37").Value = '  -4.15%  '

# Row 38
$ws.Range("E38").Value = '  +0.01%  '

# Row 39
$ws.Range("D39").Value = '''18.55'
$ws.Range("E39").Value = '  -4.21%  '

# Row 40
$ws.Range("E40").Value = '  -0.02%  '

# Row 41
$ws.Range("E41").Value = '  -3.29%  '

# Row 42
$ws.Range("E42").Value = '  -3.01%  '

# Row 43
$ws.Range("E43").Value = '  -7.49%  '

# Row 44
$ws.Range("D44").Value = '''1.18'
$ws.Range("E44").Value = '  -12.93%  '

# Row 45
$ws.Range("D45").Value = '''2.33'
$ws.Range("E45").Value = '  -9.24%  '

# Row 46
$ws.Range("D46").Value = '''38.20'
$ws.Range("E46").Value = '  -2.48%  '

# Row 47
$ws.Range("D47").Value = '''144.64'
$ws.Range("E47").Value = '  -6.16%  '

# Row 48
$ws.Range("B48").Value = 'ARBITRUM'
$ws.Range("C48").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D48").Value = '''0.531'
$ws.Range("E48").Value = '  -3.68%  '

# Row 49
$ws.Range("B49").Value = 'Filecoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D49").Value = '''3.54'
$ws.Range("E49").Value = '  -3.54%  '

# Row 50
$ws.Range("E50").Value = '  -5.32%  '

# Row 51
$ws.Range("E51").Value = '  -2.42%  '
